# Applies the "make bound-capped uncertainty distributions triangular;
# make sorbate_c,d and sorbate_g,h prod caps same as rest" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 21 & 31: Uniform -> Triangular (Shape column only) ---
$ws.Range("F21").Value = "Triangular"
$ws.Range("F31").Value = "Triangular"

# --- Row 38: Hydrogenation TAL-to-HMP conversion ---
# Shape Uniform -> Triangular; Upper-bound formula -> symmetric 0.9*E style;
# new Midpoint (H) formula referencing baseline; statement text rename.
$ws.Range("F38").Value = "Triangular"
$ws.Range("G38").Formula = "=0.9*E38"
$ws.Range("H38").Formula = "=E38"
$ws.Range("K38").Value = "R401.TAL_to_HMP_rxn.X = x"

# --- Row 40: label rename (Dehydration -> Etherification & hydrolysis) ---
$ws.Range("A40").Value = "Etherification & hydrolysis catalyst Amberlyst70:HMP ratio"

# --- Row 41: label rename ---
$ws.Range("A41").Value = "Etherification & hydrolysis reaction time"

# --- Row 42: label rename ---
$ws.Range("A42").Value = "Etherification & hydrolysis temperature"

# --- Row 43: label + statement rename; Shape -> Triangular; add Midpoint formula ---
$ws.Range("A43").Value = "Etherification & hydrolysis HMP-to-PSA conversion"
$ws.Range("F43").Value = "Triangular"
$ws.Range("H43").Formula = "=E43"
$ws.Range("K43").Value = "R402.HMP_to_PSA_rxn.X = x"

# --- Row 44: label rename ---
$ws.Range("A44").Value = "Etherification & hydrolysis pressure"

# --- Row 45: label rename ---
$ws.Range("A45").Value = "Etherification & hydrolysis spent catalyst Amberlyst70 replacement rate"

# --- Row 47: Shape -> Triangular; Upper-bound value -> formula; add Midpoint formula ---
$ws.Range("F47").Value = "Triangular"
$ws.Range("G47").Formula = "=0.9*E47"
$ws.Range("H47").Formula = "=E47"

# --- Sheet view: selection moved from A11 (whole row) to G50 ---
$ws.Range("G50").Select()
